# Update "想去人数" (F column) figures on sheets "展览" and "全部类型"
# to reflect freshly generated data (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 7154
$wsExpo.Range("F13").Value = 10
$wsExpo.Range("F14").Value = 457
$wsExpo.Range("F19").Value = 3720
$wsExpo.Range("F20").Value = 28
$wsExpo.Range("F21").Value = 248
$wsExpo.Range("F25").Value = 2368
$wsExpo.Range("F27").Value = 286
$wsExpo.Range("F34").Value = 23
$wsExpo.Range("F36").Value = 1394
$wsExpo.Range("F37").Value = 130

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 7154
$wsAll.Range("F14").Value = 10
$wsAll.Range("F15").Value = 457
$wsAll.Range("F20").Value = 3720
$wsAll.Range("F21").Value = 28
$wsAll.Range("F22").Value = 248
$wsAll.Range("F26").Value = 2368
$wsAll.Range("F28").Value = 286
$wsAll.Range("F35").Value = 23
$wsAll.Range("F37").Value = 1394
$wsAll.Range("F38").Value = 130
